$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E47").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema na Internet`n* Sucesso na Consulta de Dados`n* Navegação Reduzida`n* Sucesso no billing`nQuer pacote adicional`nEncaminha para o ATH."
$ws.Range("F47").Value = "1 - Chamada é encaminhada para a URA Cognitiva;`n2 - A URA Cognitiva informa que existem faturas em aberto e pergunta se ele tem interesse no processo de Religa;`n3 - A URA Cognitiva pergunta se o usuário está com dificuldade para usar os serviços;`n4 - A URA Cognitiva pergunta em qual serviço o usuário está com dificuldades/problemas;`n5 - A URA Cognitiva informa que a velocidade está reduzida e a data de renovação;`n6 - URA informa que vai tranferir para o ATH."
$ws.Range("E49").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema na Internet`n* Sucesso na Consulta de Dados`n* Navegação Reduzida`n* Falha no billing`nQuer pacote adicional`nEncaminha para o ATH."
$ws.Range("F49").Value = "1 - Chamada é encaminhada para a URA Cognitiva;`n2 - A URA Cognitiva informa que existem faturas em aberto e pergunta se ele tem interesse no processo de Religa;`n3 - A URA Cognitiva pergunta se o usuário está com dificuldade para usar os serviços;`n4 - A URA Cognitiva pergunta em qual serviço o usuário está com dificuldades/problemas;`n5 - A URA Cognitiva pergunta se o usuário quer saber mais sobre os pacotes adicionais;`n6 - URA informa que vai tranferir para o ATH."
$ws.Range("E51").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema na Internet`n* Sucesso na Consulta de Dados`n* Navegação Normal`nQuer falar com ATH`nEncaminha para o ATH."
$ws.Range("E52").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema na Internet`n* Sucesso na Consulta de Dados`n* Navegação Normal`nQuer continuar na URA`nJá reiniciou`n* Sucesso no Envio do Guia`nNão quer mais nada"
$ws.Range("E53").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema na Internet`n* Sucesso na Consulta de Dados`n* Navegação Normal`nQuer continuar na URA`nJá reiniciou`n* Falha no Envio do Guia`nNão quer mais nada"
$ws.Range("E54").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema na Internet`n* Sucesso na Consulta de Dados`n* Navegação Normal`nQuer continuar na URA`nNão reiniciou`nAceita Reiniciar"
$ws.Range("E55").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema na Internet`n* Sucesso na Consulta de Dados`n* Navegação Normal`nQuer continuar na URA`nNão reiniciou`nNão quer reiniciar`n* Sucesso no Envio do Guia`nNão quer mais nada"
$ws.Range("E56").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema na Internet`n* Sucesso na Consulta de Dados`n* Navegação Normal`nQuer continuar na URA`nNão reiniciou`nNão quer reiniciar`n* Falha no Envio do Guia`nNão quer mais nada"
$ws.Range("E57").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema na Internet`n* Falha na Consulta de Dados`nJá reiniciou`n* Sucesso no Envio do Guia`nNão quer mais nada"
$ws.Range("E58").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema na Internet`n* Falha na Consulta de Dados`nJá reiniciou`n* Falha no Envio do Guia`nNão quer mais nada"
$ws.Range("E59").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema na Internet`n* Falha na Consulta de Dados`nNão reiniciou`nAceita Reiniciar"
$ws.Range("E60").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema na Internet`n* Falha na Consulta de Dados`nNão reiniciou`nNão quer reiniciar`n* Sucesso no Envio do Guia`nNão quer mais nada"
$ws.Range("E61").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema na Internet`n* Falha na Consulta de Dados`nNão reiniciou`nNão quer reiniciar`n* Falha no Envio do Guia`nNão quer mais nada"
$ws.Range("E62").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema nas ligações`nQuer falar com ATH`nEncaminha para o ATH."
$ws.Range("E63").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema nas ligações`nQuer continuar na URA`nJá reiniciou`n* Sucesso no Envio do Guia`nNão quer mais nada"
$ws.Range("E64").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema nas ligações`nQuer continuar na URA`nJá reiniciou`n* Falha no Envio do Guia`nNão quer mais nada"
$ws.Range("E65").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema nas ligações`nQuer continuar na URA`nNão reiniciou`nAceita Reiniciar"
$ws.Range("E66").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema nas ligações`nQuer continuar na URA`nNão reiniciou`nNão quer reiniciar`n* Sucesso no Envio do Guia`nNão quer mais nada"
$ws.Range("E67").Value = "Quer fazer Religa`nEstá enfrentando dificuldade`nProblema nas ligações`nQuer continuar na URA`nNão reiniciou`nNão quer reiniciar`n* Falha no Envio do Guia`nNão quer mais nada"
